$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 354, pushing existing row 354 (and below) down to 355.
$ws.Rows.Item(354).Insert()

# Populate the newly inserted row 354 with the new data record.
$ws.Cells.Item(354, 1).Value = 10
$ws.Cells.Item(354, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(354, 3).Value = "La Araucanía"
$ws.Cells.Item(354, 4).Value = 44841
$ws.Cells.Item(354, 5).Value = 9
$ws.Cells.Item(354, 6).Value = 100112040
$ws.Cells.Item(354, 7).Value = "Cilantro"
$ws.Cells.Item(354, 8).Value = "Sin especificar"
$ws.Cells.Item(354, 9).Value = "Primera"
$ws.Cells.Item(354, 10).Value = 50
$ws.Cells.Item(354, 11).Value = 5000
$ws.Cells.Item(354, 12).Value = 5000
$ws.Cells.Item(354, 13).Value = 5000
$ws.Cells.Item(354, 14).Value = "`$/docena de atados (2 kilos)"
$ws.Cells.Item(354, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(354, 16).Value = 2500
$ws.Cells.Item(354, 17).Value = 2
$ws.Cells.Item(354, 18).Value = "Hortaliza"
